$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = @"
questions = [
    {
        "title": "Which of the following instructions will define a new String variable called firstvariable?",
        "ques_type": 2,
        "options": [
            "Dim firstvariable as String",
            "Set firstvariable as String",
            "var firstvariable (String)",
            "Declare firstvariable (String)"
        ],
        "score": "Dim firstvariable as String"
    },
    {
        "title": "What result should you expect if you apply the instruction Activecell.Value = 5 with the selection of cells shown in the screenshot below?",
        "ques_type": 2,
        "options": [
            "Make the values of all of the cells in the range A4:C8 equal to 5.",
            "Add 5 to the existing values of all the cells in the range A4:C8.",
            "Make the value of cell A4 equal to 5.",
            "Add 5 to the value of cell A4."
        ],
        "score": "Make the value of cell A4 equal to 5."
    },
    {
        "title": "True or false: New labels on UserForms can be created at runtime.",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    },
    {
        "title": "The following scripts (ScriptA and ScriptB, included below) make the same results. True or false: ScriptB performs faster than ScriptA. Sub ScriptA()\nWhile ActiveCell.Value = \"\"\n           ActiveCell.Value = 1\n           ActiveCell.Offset(1, 0).Select\nWend\nEnd Sub\n \nSub ScriptB()\ni = 0\nDo\n           ActiveCell.Offset(i, 0).Value = 1\n           i = i + 1 \nLoop While ActiveCell.Offset(i, 0).Value = \"\"\nEnd Sub",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "True"
    }
]
"@

# Remove trailing newline that here-string syntax may add before the closing "@
if ($text.EndsWith("`r`n")) {
    $text = $text.Substring(0, $text.Length - 2)
} elseif ($text.EndsWith("`n")) {
    $text = $text.Substring(0, $text.Length - 1)
}

# Clear old formatting (remove bold font / border from A1) and reset A2
$ws.Range("A1:A2").ClearContents()
$ws.Range("A1:A2").ClearFormats()

$ws.Range("A1").Value = $text
